# Updated symbol list on Wed Dec 21 21:29:08 UTC 2022 with GitHub Actions
# Update cryptocurrency price (column D) values to the latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.87"
$ws.Range("D4").Value = "'5.238"
$ws.Range("D5").Value = "'0.05695"
$ws.Range("D6").Value = "'3.413"
$ws.Range("D7").Value = "'6.307"
$ws.Range("D8").Value = "'0.8062"
$ws.Range("D9").Value = "'0.8590"
$ws.Range("D10").Value = "'0.1412"
$ws.Range("D11").Value = "'0.07345"
$ws.Range("D12").Value = "'0.03030"
$ws.Range("D13").Value = "'0.03094"
$ws.Range("D15").Value = "'3.868"
$ws.Range("D16").Value = "'0.001592"
$ws.Range("D17").Value = "'0.04766"
$ws.Range("D18").Value = "'0.0005839"
$ws.Range("D20").Value = "'0.005028"
$ws.Range("D21").Value = "'0.0009962"
$ws.Range("D22").Value = "'0.0001499"
$ws.Range("D24").Value = "'2.193"
$ws.Range("D25").Value = "'0.3280"
$ws.Range("D26").Value = "'0.1282"
$ws.Range("D40").Value = "'0.03915"
$ws.Range("D41").Value = "'0.006841"
$ws.Range("D42").Value = "'0.1066"
$ws.Range("D43").Value = "'0.003199"
$ws.Range("D44").Value = "'0.008464"
$ws.Range("D45").Value = "'0.00005594"
$ws.Range("D47").Value = "'0.4500"
$ws.Range("D48").Value = "'0.1991"
